$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 55.77111053466797
$ws.Range("E2").Value = 55.10889053344727
$ws.Range("F2").Value = 62.11555480957031
$ws.Range("G2").Value = 54.61000061035156
$ws.Range("H2").Value = 354496073
$ws.Range("I2").Value = "ISRG"

$ws.Range("D3").Value = 54.15888977050781
$ws.Range("E3").Value = 59.24111175537109
$ws.Range("F3").Value = 62.76222229003906
$ws.Range("G3").Value = 52.31000137329102
$ws.Range("H3").Value = 354496073
$ws.Range("I3").Value = "ISRG"

$ws.Range("D4").Value = 51.11000061035156
$ws.Range("E4").Value = 55.17777633666992
$ws.Range("F4").Value = 57.77555465698242
$ws.Range("G4").Value = 49.66666793823242
$ws.Range("H4").Value = 354496073
$ws.Range("I4").Value = "ISRG"

$ws.Range("D5").Value = 60.79555511474609
$ws.Range("E5").Value = 60.09444427490234
$ws.Range("F5").Value = 64.44444274902344
$ws.Range("G5").Value = 57.83333206176758
$ws.Range("H5").Value = 354496073
$ws.Range("I5").Value = "ISRG"

$ws.Range("D6").Value = 66.13222503662109
$ws.Range("E6").Value = 69.59555816650391
$ws.Range("F6").Value = 72.76444244384766
$ws.Range("G6").Value = 65.67111206054688
$ws.Range("H6").Value = 354496073
$ws.Range("I6").Value = "ISRG"

$ws.Range("D7").Value = 73.50777435302734
$ws.Range("E7").Value = 77.30666351318359
$ws.Range("F7").Value = 80.00111389160156
$ws.Range("G7").Value = 73.22555541992188
$ws.Range("H7").Value = 354496073
$ws.Range("I7").Value = "ISRG"

$ws.Range("D8").Value = 80.17333221435547
$ws.Range("E8").Value = 74.67555236816406
$ws.Range("F8").Value = 80.80555725097656
$ws.Range("G8").Value = 73.08000183105469
$ws.Range("H8").Value = 354496073
$ws.Range("I8").Value = "ISRG"

$ws.Range("D9").Value = 70.66000366210938
$ws.Range("E9").Value = 76.96555328369141
$ws.Range("F9").Value = 77.82444763183594
$ws.Range("G9").Value = 69.413330078125
$ws.Range("H9").Value = 354496073
$ws.Range("I9").Value = "ISRG"

$ws.Range("D10").Value = 85.07111358642578
$ws.Range("E10").Value = 92.87444305419922
$ws.Range("F10").Value = 93.82888793945312
$ws.Range("G10").Value = 83.70888519287109
$ws.Range("H10").Value = 354496073
$ws.Range("I10").Value = "ISRG"

$ws.Range("D11").Value = 104.5511093139648
$ws.Range("E11").Value = 104.2511138916016
$ws.Range("F11").Value = 108.2955551147461
$ws.Range("G11").Value = 102.2855529785156
$ws.Range("H11").Value = 354496073
$ws.Range("I11").Value = "ISRG"

$ws.Range("D12").Value = 116.5666656494141
$ws.Range("E12").Value = 125.120002746582
$ws.Range("F12").Value = 127.3833312988281
$ws.Range("G12").Value = 116.3966674804688
$ws.Range("H12").Value = 354496073
$ws.Range("I12").Value = "ISRG"

$ws.Range("D13").Value = 122.629997253418
$ws.Range("E13").Value = 143.8899993896484
$ws.Range("F13").Value = 150.6666717529297
$ws.Range("G13").Value = 121.7033309936523
$ws.Range("H13").Value = 354496073
$ws.Range("I13").Value = "ISRG"

$ws.Range("D14").Value = 137.4199981689453
$ws.Range("E14").Value = 146.9266662597656
$ws.Range("F14").Value = 157.9299926757812
$ws.Range("G14").Value = 131.1900024414062
$ws.Range("H14").Value = 354496073
$ws.Range("I14").Value = "ISRG"

$ws.Range("D15").Value = 158.82666015625
$ws.Range("E15").Value = 169.3966674804688
$ws.Range("F15").Value = 179.7666625976562
$ws.Range("G15").Value = 157.2666625976562
$ws.Range("H15").Value = 354496073
$ws.Range("I15").Value = "ISRG"

$ws.Range("D16").Value = 191.7166595458984
$ws.Range("E16").Value = 173.7266693115234
$ws.Range("F16").Value = 193.7066650390625
$ws.Range("G16").Value = 155.6900024414062
$ws.Range("H16").Value = 354496073
$ws.Range("I16").Value = "ISRG"

$ws.Range("D17").Value = 156.5
$ws.Range("E17").Value = 174.5466613769531
$ws.Range("F17").Value = 181.25
$ws.Range("G17").Value = 148.8200073242188
$ws.Range("H17").Value = 354496073
$ws.Range("I17").Value = "ISRG"

$ws.Range("D18").Value = 191.6666717529297
$ws.Range("E18").Value = 170.2100067138672
$ws.Range("F18").Value = 196.4400024414062
$ws.Range("G18").Value = 162.2200012207031
$ws.Range("H18").Value = 354496073
$ws.Range("I18").Value = "ISRG"

$ws.Range("D19").Value = 177.3333282470703
$ws.Range("E19").Value = 173.1699981689453
$ws.Range("F19").Value = 181.6166687011719
$ws.Range("G19").Value = 171.4199981689453
$ws.Range("H19").Value = 354496073
$ws.Range("I19").Value = "ISRG"

$ws.Range("D20").Value = 179.9799957275391
$ws.Range("E20").Value = 184.3166656494141
$ws.Range("F20").Value = 189.9933319091797
$ws.Range("G20").Value = 167.336669921875
$ws.Range("H20").Value = 354496073
$ws.Range("I20").Value = "ISRG"

$ws.Range("D21").Value = 198.5666656494141
$ws.Range("E21").Value = 186.5933380126953
$ws.Range("F21").Value = 205.5200042724609
$ws.Range("G21").Value = 185.009994506836
$ws.Range("H21").Value = 354496073
$ws.Range("I21").Value = "ISRG"

$ws.Range("D22").Value = 153.4199981689453
$ws.Range("E22").Value = 170.2933349609375
$ws.Range("F22").Value = 176.5
$ws.Range("G22").Value = 150
$ws.Range("H22").Value = 354496073
$ws.Range("I22").Value = "ISRG"

$ws.Range("D23").Value = 190.6666717529297
$ws.Range("E23").Value = 228.479995727539
$ws.Range("F23").Value = 234.6666717529297
$ws.Range("G23").Value = 188.7366638183593
$ws.Range("H23").Value = 354496073
$ws.Range("I23").Value = "ISRG"

$ws.Range("D24").Value = 237.7033386230469
$ws.Range("E24").Value = 222.3600006103516
$ws.Range("F24").Value = 257.2066650390625
$ws.Range("G24").Value = 217.6666717529297
$ws.Range("H24").Value = 354496073
$ws.Range("I24").Value = "ISRG"

$ws.Range("D25").Value = 275.2799987792969
$ws.Range("E25").Value = 249.2133331298828
$ws.Range("F25").Value = 275.2799987792969
$ws.Range("G25").Value = 245.2133331298828
$ws.Range("H25").Value = 354496073
$ws.Range("I25").Value = "ISRG"

$ws.Range("D26").Value = 249.9966735839844
$ws.Range("E26").Value = 288.3333435058594
$ws.Range("F26").Value = 297.9299926757812
$ws.Range("G26").Value = 248.6266632080078
$ws.Range("H26").Value = 354496073
$ws.Range("I26").Value = "ISRG"

$ws.Range("D27").Value = 305.8433227539062
$ws.Range("E27").Value = 330.4866638183594
$ws.Range("F27").Value = 331.6966552734375
$ws.Range("G27").Value = 305.4966735839844
$ws.Range("H27").Value = 354496073
$ws.Range("I27").Value = "ISRG"

$ws.Range("D28").Value = 331.3333435058594
$ws.Range("E28").Value = 361.1300048828125
$ws.Range("F28").Value = 362.6700134277344
$ws.Range("G28").Value = 321.0533447265625
$ws.Range("H28").Value = 354496073
$ws.Range("I28").Value = "ISRG"

$ws.Range("D29").Value = 358.6499938964844
$ws.Range("E29").Value = 284.1799926757812
$ws.Range("F29").Value = 362
$ws.Range("G29").Value = 254.1999969482422
$ws.Range("H29").Value = 354496073
$ws.Range("I29").Value = "ISRG"

$ws.Range("D30").Value = 304
$ws.Range("E30").Value = 239.3000030517578
$ws.Range("F30").Value = 308.9700012207031
$ws.Range("G30").Value = 235.0700073242188
$ws.Range("H30").Value = 354496073
$ws.Range("I30").Value = "ISRG"

$ws.Range("D31").Value = 201.0700073242188
$ws.Range("E31").Value = 230.1699981689453
$ws.Range("F31").Value = 231.259994506836
$ws.Range("G31").Value = 196.1000061035156
$ws.Range("H31").Value = 354496073
$ws.Range("I31").Value = "ISRG"

$ws.Range("D32").Value = 189.3999938964844
$ws.Range("E32").Value = 246.4700012207031
$ws.Range("F32").Value = 249.1399993896484
$ws.Range("G32").Value = 180.0700073242188
$ws.Range("H32").Value = 354496073
$ws.Range("I32").Value = "ISRG"

$ws.Range("D33").Value = 269.5899963378906
$ws.Range("E33").Value = 245.6900024414062
$ws.Range("F33").Value = 273.5799865722656
$ws.Range("G33").Value = 238.8399963378907
$ws.Range("H33").Value = 354496073
$ws.Range("I33").Value = "ISRG"

$ws.Range("D34").Value = 254.3000030517578
$ws.Range("E34").Value = 301.2200012207031
$ws.Range("F34").Value = 304.8399963378906
$ws.Range("G34").Value = 250.009994506836
$ws.Range("H34").Value = 354496073
$ws.Range("I34").Value = "ISRG"

$ws.Range("D35").Value = 339.9200134277344
$ws.Range("E35").Value = 324.3999938964844
$ws.Range("F35").Value = 358.0700073242188
$ws.Range("G35").Value = 321.9100036621094
$ws.Range("H35").Value = 354496073
$ws.Range("I35").Value = "ISRG"

$ws.Range("D36").Value = 291.9599914550781
$ws.Range("E36").Value = 262.2200012207031
$ws.Range("F36").Value = 301.0700073242188
$ws.Range("G36").Value = 254.8500061035156
$ws.Range("H36").Value = 354496073
$ws.Range("I36").Value = "ISRG"

$ws.Range("D37").Value = 333
$ws.Range("E37").Value = 378.2200012207031
$ws.Range("F37").Value = 384
$ws.Range("G37").Value = 320.260009765625
$ws.Range("H37").Value = 354496073
$ws.Range("I37").Value = "ISRG"

$ws.Range("D38").Value = 399
$ws.Range("E38").Value = 370.6199951171875
$ws.Range("F38").Value = 399.6400146484375
$ws.Range("G38").Value = 364.1700134277344
$ws.Range("H38").Value = 354496073
$ws.Range("I38").Value = "ISRG"

$ws.Range("D39").Value = 445.8500061035156
$ws.Range("E39").Value = 444.6099853515625
$ws.Range("F39").Value = 468.7799987792969
$ws.Range("G39").Value = 413.8200073242188
$ws.Range("H39").Value = 354496073
$ws.Range("I39").Value = "ISRG"

$ws.Range("D40").Value = 492.5
$ws.Range("E40").Value = 503.8399963378906
$ws.Range("F40").Value = 523.3400268554688
$ws.Range("G40").Value = 470.0599975585938
$ws.Range("H40").Value = 354496073
$ws.Range("I40").Value = "ISRG"

$ws.Range("D41").Value = 521.1900024414062
$ws.Range("E41").Value = 571.8800048828125
$ws.Range("F41").Value = 616
$ws.Range("G41").Value = 520.9199829101562
$ws.Range("H41").Value = 354496073
$ws.Range("I41").Value = "ISRG"

$ws.Range("D42").Value = 492.260009765625
$ws.Range("E42").Value = 515.7999877929688
$ws.Range("F42").Value = 529.1900024414062
$ws.Range("G42").Value = 425
$ws.Range("H42").Value = 354496073
$ws.Range("I42").Value = "ISRG"

$ws.Range("D43").Value = 542.5800170898438
$ws.Range("E43").Value = 481.0899963378906
$ws.Range("F43").Value = 550.7000122070312
$ws.Range("G43").Value = 479.6000061035156
$ws.Range("H43").Value = 354496073
$ws.Range("I43").Value = "ISRG"

$ws.Range("D44").Value = 447
$ws.Range("E44").Value = 546.510009765625
$ws.Range("F44").Value = 552.5
$ws.Range("G44").Value = 427.8399963378906
$ws.Range("H44").Value = 354496073
$ws.Range("I44").Value = "ISRG"
